$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.724.82'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '3.087.11'
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''517.51'
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").Value = '''143.64'
$ws.Range("E6").Value = '  +4.14%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").Value = '''7.34'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("E10").Value = '  +0.37%  '
$ws.Range("D11").Value = '''0.375'
$ws.Range("E11").Value = '  +2.21%  '
$ws.Range("D12").Value = '3.605.73'
$ws.Range("E12").Value = '  +2.91%  '
$ws.Range("E13").Value = '  +2.62%  '
$ws.Range("D14").Value = '''25.81'
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").Value = '57.772.23'
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D17").Value = '''6.17'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '3.076.53'
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("D19").Value = '''13.10'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").Value = '''8.22'
$ws.Range("E20").Value = '  +2.21%  '
$ws.Range("D21").Value = '''337.81'
$ws.Range("E21").Value = '  +3.48%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").Value = '''0.502'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '''65.69'
$ws.Range("E24").Value = '  +2.39%  '
$ws.Range("E25").Value = '  +5.77%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '0.0₃0933'
$ws.Range("E27").Value = '  +5.99%  '
$ws.Range("D28").Value = '''6.49'
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("D29").Value = '''7.11'
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").Value = '''1.82'
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").Value = '''20.88'
$ws.Range("E31").Value = '  +1.21%  '
$ws.Range("E32").Value = '  -2.32%  '
$ws.Range("D33").Value = '''154.47'
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("D34").Value = '''4.54'
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("D35").Value = '''5.93'
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("D36").Value = '''26.74'
$ws.Range("E36").Value = '  +5.98%  '
$ws.Range("E37").Value = '  +1.06%  '
$ws.Range("D38").Value = '''0.0689'
$ws.Range("E38").Value = '  +3.18%  '
$ws.Range("D39").Value = '3.124.27'
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("D40").Value = '''36.95'
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").Value = '''3.88'
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("D42").Value = '''0.673'
$ws.Range("E42").Value = '  +3.40%  '
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").Value = '2.278.18'
$ws.Range("E44").Value = '  +5.00%  '
$ws.Range("E45").Value = '  +5.03%  '
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D49").Value = '''5.89'
$ws.Range("E49").Value = '  -4.07%  '
$ws.Range("D50").Value = '''0.0877'
$ws.Range("E50").Value = '  +2.74%  '
$ws.Range("D51").Value = '''0.693'
$ws.Range("E51").Value = '  +3.43%  '

# Row 47/48 swap (ONDO <-> InjectiveProtocol) with updated data
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '''20.46'
$ws.Range("E47").Value = '  +4.83%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '''0.952'
$ws.Range("E48").Value = '  +0.87%  '
